# Append a new data row (row 14) to the active worksheet, mirroring the
# existing rows of stock-prediction "Bag" data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 14

$ws.Cells.Item($row, 1).Value = 42620.889513888891
$ws.Cells.Item($row, 2).Value = 14
$ws.Cells.Item($row, 3).Value = 59
$ws.Cells.Item($row, 4).Value = 38
$ws.Cells.Item($row, 5).Value = 59
$ws.Cells.Item($row, 6).Value = 16
$ws.Cells.Item($row, 7).Value = 21184
$ws.Cells.Item($row, 8).Value = 19161
$ws.Cells.Item($row, 9).Value = 2172
$ws.Cells.Item($row, 10).Value = 272
$ws.Cells.Item($row, 11).Value = 173
$ws.Cells.Item($row, 12).Value = 20
$ws.Cells.Item($row, 13).Value = 4
$ws.Cells.Item($row, 14).Value = "Bag"

# Match the date/time number formatting used by the other rows in column A
$ws.Cells.Item($row, 1).NumberFormat = "m/d/yy h:mm"
